$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values per row (2-10) for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$data = @{
    2  = @{ E=3; G=7.636417666666667; H=22.909253; I=0.108532481296676; J=0.108532481296676; K=3; M=45.95651366666667; N=137.869541; O=0.6189188856627118; P=0.6189188856627118; Q=350.9431328625415; R=3158.488195762873; S=0.06717280238234781; T=0.06717280238234781 }
    3  = @{ E=3; G=7.636417666666667; H=22.909253; I=0.108532481296676; J=0.108532481296676; K=3; M=6.849914666666667; N=20.549744; O=0.09225115688993263; P=0.09225115688993261; Q=52.30880937569245; R=470.779284381232; S=0.01001224695975333; T=0.01001224695975333 }
    4  = @{ E=3; G=7.636417666666667; H=22.909253; I=0.108532481296676; J=0.108532481296676; K=3; M=21.446458; N=64.33937399999999; O=0.2888299574473556; P=0.2888299574473556; Q=163.7741107586246; R=1473.966996827622; S=0.03134743195457484; T=0.03134743195457484 }
    5  = @{ E=3; G=15.103385; H=45.31015499999999; I=0.214656652056136; J=0.214656652056136; K=3; M=45.95651366666667; N=137.869541; O=0.6189188856627118; P=0.6189188856627118; Q=694.0989191654282; R=6246.890272488854; S=0.1328550558906721; T=0.1328550558906721 }
    6  = @{ E=3; G=15.103385; H=45.31015499999999; I=0.214656652056136; J=0.214656652056136; K=3; M=6.849914666666667; N=20.549744; O=0.09225115688993263; P=0.09225115688993261; Q=103.4568984278133; R=931.1120858503199; S=0.01980232448629828; T=0.01980232448629828 }
    7  = @{ E=3; G=15.103385; H=45.31015499999999; I=0.214656652056136; J=0.214656652056136; K=3; M=21.446458; N=64.33937399999999; O=0.2888299574473556; P=0.2888299574473556; Q=323.9141120603299; R=2915.227008542969; S=0.06199927167916557; T=0.06199927167916557 }
    8  = @{ E=3; G=47.62086333333334; H=142.86259; I=0.676810866647188; J=0.676810866647188; K=3; M=45.95651366666667; N=137.869541; O=0.6189188856627118; P=0.6189188856627118; Q=2188.488856596799; R=19696.39970937119; S=0.4188910273896919; T=0.4188910273896919 }
    9  = @{ E=3; G=47.62086333333334; H=142.86259; I=0.676810866647188; J=0.676810866647188; K=3; M=6.849914666666667; N=20.549744; O=0.09225115688993263; P=0.09225115688993261; Q=326.1988501863289; R=2935.78965167696; S=0.06243658544388101; T=0.062436585443881 }
    10 = @{ E=3; G=47.62086333333334; H=142.86259; I=0.676810866647188; J=0.676810866647188; K=3; M=21.446458; N=64.33937399999999; O=0.2888299574473556; P=0.2888299574473556; Q=1021.298845402073; R=9191.689608618659; S=0.1954832538136152; T=0.1954832538136152 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
